# Add new negative/positive test-data columns to the "CreateUser" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateUser")

# New columns, added one at a time (header then value) to mirror the
# original authoring order of the shared-string table.
$ws.Range("G1").Value = "ShortPassword"
$ws.Range("G2").Value = "test12"

$ws.Range("H1").Value = "PasswordContaintUserName"
$ws.Range("H2").Value = "Mart123456"

$ws.Range("I1").Value = "WrongConfirmedPassword"
$ws.Range("I2").Value = "test1234567"

# Move the active selection to match the source change (I11)
$ws.Range("I11").Select()
